# Generate Report for Handoff
#
# This updates the localization-status report for the file
# "3e883624-a006-461c-976c-f3d09aed599d.md" with a freshly generated
# handoff xliff / report timestamp:
#   - Overview sheet: "Latest HO Xliff Generate Date" -> 2016-09-07 04:55:09
#   - zh-cn sheet:   "Latest Handoff Datetime"        -> 2016-09-07 04:54:57
#   - de-de sheet:   "Latest Handoff Datetime"        -> 2016-09-07 04:55:09

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Row 6 on every sheet corresponds to 3e883624-a006-461c-976c-f3d09aed599d.md
$wsOverview.Range("G6").Value = "2016-09-07 04:55:09"
$wsZhCn.Range("H6").Value     = "2016-09-07 04:54:57"
$wsDeDe.Range("H6").Value     = "2016-09-07 04:55:09"
